$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("Q4").Value = "0000-1001"
$ws.Range("Q5").Value = "0000-1002"
$ws.Range("Q6").Value = "0000-1003"
$ws.Range("Q7").Value = "0000-1004"
$ws.Range("Q8").Value = "0000-1005"
$ws.Range("Q9").Value = "0000-1006"

$ws.Application.ActiveWindow.Panes.Item(4).ScrollColumn = 10
$ws.Range("Q10").Select()
